$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.012.45'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '2.932.19'
$ws.Range("E3").Value = '  +4.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '203.28'
$ws.Range("E5").Value = '  +8.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '597.72'
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.198'
$ws.Range("E9").Value = '  +5.14%  '
$ws.Range("D10").Value = '2.928.30'
$ws.Range("E10").Value = '  +4.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.437'
$ws.Range("E11").Value = '  +17.18%  '
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("D14").Value = '3.468.15'
$ws.Range("E14").Value = '  +4.07%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '75.879.80'
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.08'
$ws.Range("E16").Value = '  +5.04%  '
$ws.Range("E17").Value = '  +2.73%  '
$ws.Range("D18").Value = '2.927.38'
$ws.Range("E18").Value = '  +4.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.19'
$ws.Range("E19").Value = '  +7.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.88'
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.49'
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("E22").Value = '  +2.49%  '
$ws.Range("E23").Value = '  +5.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.72'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  +4.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.30'
$ws.Range("E27").Value = '  +4.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.68'
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("E29").Value = '  +6.04%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '503.76'
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.79'
$ws.Range("E33").Value = '  +2.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  +3.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.28'
$ws.Range("E36").Value = '  +2.12%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.85'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("B38").Value = 'Cronos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.111'
$ws.Range("E38").Value = '  +29.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.63'
$ws.Range("E39").Value = '  +1.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.371'
$ws.Range("E40").Value = '  +9.66%  '
$ws.Range("E41").Value = '  -4.10%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '181.47'
$ws.Range("E43").Value = '  -0.69%  '
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.17'
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("E47").Value = '  -0.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.35'
$ws.Range("E48").Value = '  +1.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.581'
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("E50").Value = '  +1.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.45'
$ws.Range("E51").Value = '  +8.18%  '
